# Ind_Customer_Month.xlsx - update title to reflect new month (March 2025)
# and refresh the hourly demand (kWh) values for the 24 hours.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A1").Value = "An Individual Customer's Average Daily Load Profile Usage for March 2025"

$values = @{
    3  = 0.0007346346064516128
    4  = 0.000508684864516129
    5  = 0.0002828067483870968
    6  = 0.0002506920129032259
    7  = 0.0001861042193548387
    8  = 0.0002183622838709678
    9  = 0.0002506203483870968
    10 = 0.0002506203483870968
    11 = 0.0003870967741935484
    12 = 0.0003870967741935484
    13 = 0.0003870967741935484
    14 = 0.0005161290322580645
    15 = 0.0006451612903225806
    16 = 0.0006129032258064517
    17 = 0.0006450717096774194
    18 = 0.0006452508709677419
    19 = 0.0004193010903225807
    20 = 0.0003870878161290323
    21 = 0.0005806272225806452
    22 = 0.0007097580645161291
    23 = 0.0007096774193548386
    24 = 0.0006129032258064516
    25 = 0.0005806451612903226
    26 = 0.000709534129032258
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 2).Value = $values[$row]
}
